# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text ("Ready for handoff" -> "Handed back: in sync with en-US") is
#    updated on the Overview summary columns and on each locale detail sheet
#  - zh-cn / de-de sheets get their "Latest Target File", "Latest Handback File"
#    and "Latest Handback DateTime" columns populated with the real handback data
#  - A hyperlink to the source markdown file is added on the "Latest Target File" cell

$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/220eccca9608598735460c3e731317c4792866d4/e2e/35c6a7cd-6bef-496f-888f-e35934852f7c.md"
$mdDisplay = "35c6a7cd-6bef-496f-888f-e35934852f7c.md"

# ---------------------------------------------------------------------------
# Overview sheet: status text for zh-cn / de-de moves from "Ready for handoff"
# to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("J2").Value = "35c6a7cd-6bef-496f-888f-e35934852f7c.33998104ed7215556f29712780f9becdee53c75f.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-19 04:53:25"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("I2").Font.ThemeFont = 0

$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("J2").Value = "35c6a7cd-6bef-496f-888f-e35934852f7c.33998104ed7215556f29712780f9becdee53c75f.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-19 04:53:32"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("I2").Font.ThemeFont = 0

$wsDeDe.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17
